{"js": "// Update the worksheet date and each three-digit-by-one-digit division\n// answer (\"dividend\u00f7divisor=quotient, remainder\") to its new value.\n// Every source string in this document is unique, so a plain\n// search-and-replace per pair is sufficient \u2014 no row/column bookkeeping\n// is needed.\nconst pairs = [\n  [\"2024-07-31 Wednesday\", \"2024-08-01 Thursday\"],\n  [\"664\u00f76=110, 4\", \"742\u00f75=148, 2\"],\n  [\"204\u00f74=51, 0\", \"540\u00f77=77, 1\"],\n  [\"181\u00f78=22, 5\", \"943\u00f77=134, 5\"],\n  [\"243\u00f78=30, 3\", \"891\u00f78=111, 3\"],\n  [\"860\u00f72=430, 0\", \"189\u00f72=94, 1\"],\n  [\"399\u00f74=99, 3\", \"720\u00f79=80, 0\"],\n  [\"882\u00f73=294, 0\", \"842\u00f79=93, 5\"],\n  [\"519\u00f76=86, 3\", \"330\u00f77=47, 1\"],\n  [\"411\u00f79=45, 6\", \"364\u00f72=182, 0\"],\n  [\"226\u00f73=75, 1\", \"126\u00f72=63, 0\"],\n  [\"616\u00f75=123, 1\", \"663\u00f72=331, 1\"],\n  [\"920\u00f74=230, 0\", \"827\u00f78=103, 3\"],\n  [\"657\u00f77=93, 6\", \"978\u00f74=244, 2\"],\n  [\"766\u00f72=383, 0\", \"984\u00f77=140, 4\"],\n  [\"131\u00f74=32, 3\", \"242\u00f76=40, 2\"],\n  [\"456\u00f79=50, 6\", \"568\u00f75=113, 3\"],\n  [\"549\u00f78=68, 5\", \"855\u00f77=122, 1\"],\n  [\"817\u00f78=102, 1\", \"317\u00f78=39, 5\"],\n  [\"227\u00f79=25, 2\", \"624\u00f75=124, 4\"],\n  [\"967\u00f76=161, 1\", \"601\u00f76=100, 1\"],\n  [\"556\u00f72=278, 0\", \"107\u00f79=11, 8\"],\n  [\"706\u00f77=100, 6\", \"237\u00f74=59, 1\"],\n  [\"912\u00f76=152, 0\", \"258\u00f79=28, 6\"],\n  [\"846\u00f79=94, 0\", \"748\u00f72=374, 0\"],\n  [\"332\u00f73=110, 2\", \"259\u00f73=86, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and each three-digit-by-one-digit division\n# answer (\"dividend\u00f7divisor=quotient, remainder\") to its new value via\n# Find/Replace. Every source string in this document is unique, so a plain\n# Execute(FindText, ..., Forward:=True, Wrap:=wdFindContinue(1),\n#         ReplaceWith:=NewText, Replace:=wdReplaceAll(2)) per pair is\n# sufficient - no paragraph/cell bookkeeping needed.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-07-31 Wednesday\", \"2024-08-01 Thursday\"),\n    @(\"664\u00f76=110, 4\", \"742\u00f75=148, 2\"),\n    @(\"204\u00f74=51, 0\", \"540\u00f77=77, 1\"),\n    @(\"181\u00f78=22, 5\", \"943\u00f77=134, 5\"),\n    @(\"243\u00f78=30, 3\", \"891\u00f78=111, 3\"),\n    @(\"860\u00f72=430, 0\", \"189\u00f72=94, 1\"),\n    @(\"399\u00f74=99, 3\", \"720\u00f79=80, 0\"),\n    @(\"882\u00f73=294, 0\", \"842\u00f79=93, 5\"),\n    @(\"519\u00f76=86, 3\", \"330\u00f77=47, 1\"),\n    @(\"411\u00f79=45, 6\", \"364\u00f72=182, 0\"),\n    @(\"226\u00f73=75, 1\", \"126\u00f72=63, 0\"),\n    @(\"616\u00f75=123, 1\", \"663\u00f72=331, 1\"),\n    @(\"920\u00f74=230, 0\", \"827\u00f78=103, 3\"),\n    @(\"657\u00f77=93, 6\", \"978\u00f74=244, 2\"),\n    @(\"766\u00f72=383, 0\", \"984\u00f77=140, 4\"),\n    @(\"131\u00f74=32, 3\", \"242\u00f76=40, 2\"),\n    @(\"456\u00f79=50, 6\", \"568\u00f75=113, 3\"),\n    @(\"549\u00f78=68, 5\", \"855\u00f77=122, 1\"),\n    @(\"817\u00f78=102, 1\", \"317\u00f78=39, 5\"),\n    @(\"227\u00f79=25, 2\", \"624\u00f75=124, 4\"),\n    @(\"967\u00f76=161, 1\", \"601\u00f76=100, 1\"),\n    @(\"556\u00f72=278, 0\", \"107\u00f79=11, 8\"),\n    @(\"706\u00f77=100, 6\", \"237\u00f74=59, 1\"),\n    @(\"912\u00f76=152, 0\", \"258\u00f79=28, 6\"),\n    @(\"846\u00f79=94, 0\", \"748\u00f72=374, 0\"),\n    @(\"332\u00f73=110, 2\", \"259\u00f73=86, 1\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #         Format, ReplaceWith, Replace)\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
